$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9599575996398926
$ws.Range("B1").Value = 1.983177900314331
$ws.Range("C1").Value = 2.891309261322021
$ws.Range("D1").Value = 3.51329779624939
$ws.Range("E1").Value = 2.060942649841309
